$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.365.04'
$ws.Range('E2').Value = '  +7.93%  '
$ws.Range('D3').Value = '2.634.42'
$ws.Range('E3').Value = '  +7.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '185.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +14.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '583.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.25%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  +4.60%  '
$ws.Range('E9').Value = '  +18.30%  '
$ws.Range('D10').Value = '2.630.23'
$ws.Range('E10').Value = '  +7.87%  '
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('E12').Value = '  +8.08%  '
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('E14').Value = '  +6.49%  '
$ws.Range('D15').Value = '74.177.50'
$ws.Range('E15').Value = '  +7.89%  '
$ws.Range('D16').Value = '3.112.87'
$ws.Range('E16').Value = '  +7.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +12.86%  '
$ws.Range('D18').Value = '2.631.79'
$ws.Range('E18').Value = '  +7.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +30.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '372.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.69%  '
$ws.Range('E22').Value = '  +18.29%  '
$ws.Range('E23').Value = '  +6.70%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.75%  '
$ws.Range('D28').Value = '2.765.47'
$ws.Range('E28').Value = '  +7.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('D30').Value = '0.0₃0950'
$ws.Range('E30').Value = '  +15.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '527.98'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +22.76%  '
$ws.Range('E32').Value = '  +19.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.76%  '
$ws.Range('E34').Value = '  +8.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.15'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.98%  '
$ws.Range('E37').Value = '  +12.84%  '
$ws.Range('E38').Value = '  +6.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.28'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +12.54%  '
$ws.Range('E42').Value = '  +9.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '160.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +23.47%  '
$ws.Range('E45').Value = '  +11.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +15.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0855'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +18.62%  '
$ws.Range('E49').Value = '  +8.85%  '
$ws.Range('E50').Value = '  +9.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +24.49%  '
